# Update results with corrected TwelveBarBreakout
# - Regenerated all 28 stock analyses; TwelveBarBreakout now checks the
#   previous swing low only.
# - Append 23 new "twelve_bar_breakout" signal instances (rows 70-92) to the
#   "Signal Instances" sheet.
# - Insert a new "twelve_bar_breakout" summary row (row 4) into the
#   "Summary" sheet, pushing the existing "random_baseline" row down to
#   row 5.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Signal Instances" sheet: append the 23 new twelve_bar_breakout rows
# ---------------------------------------------------------------------------
$wsSignals = $wb.Worksheets.Item("Signal Instances")

$newSignalRows = New-Object System.Collections.ArrayList
$null = $newSignalRows.Add(@("twelve_bar_breakout", 28, "", 1.246038794517517, 7.8675107493104, $true, 12.60892816357084, 50, -19.55717480474675, 7, $true, 5, -9.963089322576513, 4.74141741426044, -17.83060007188691, -22.57201748614735))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 52, "", 1.232639193534851, 23.16006370134525, $true, 31.39815881539322, 48, -7.335867344762787, 9, $true, 9, -7.335867344762787, 8.238095114047962, -30.49593104610804, -38.73402616015601))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 73, "", 1.363763928413391, 27.76974218976247, $true, 28.92081987684699, 52, -5.415135524404695, 2, $false, "", "", 1.15107768708452, "", ""))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 149, "", 1.687347769737244, 14.23075447601501, $true, 22.78096548929607, 44, -9.615385279482428, 8, $true, 36, 7.696125450996181, 8.55021101328106, -6.53462902501883, -15.08484003829989))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 176, "", 1.789251089096069, 25.8016915616331, $true, 32.422826369288, 45, -4.062500316109237, 9, $true, 9, 1.562501873834987, 6.621134807654908, -24.23918968779811, -30.86032449545302))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 203, "", 2.042887926101685, 23.35487698945425, $true, 25.29160991730906, 52, -4.519774375072537, 3, $true, 22, 8.442954904672245, 1.936732927854813, -14.911922084782, -16.84865501263682))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 211, "", 2.0832839012146, 21.54733285679002, $true, 22.86217590166646, 45, -4.016637602638463, 2, $true, 14, 6.340188735862469, 1.314843044876444, -15.20714412092755, -16.52198716580399))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 243, "", 2.345683336257935, 38.52537005353149, $true, 44.10889333211798, 51, -8.916803833240392, 3, $true, 3, 0.1655637048255991, 5.58352327858649, -38.35980634870589, -43.94332962729238))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 270, "", 2.647828102111816, 8.963319693608634, $true, 27.66456758507744, 24, -3.18379171797038, 40, $true, 36, 5.59580781062429, 18.7012478914688, -3.367511882984344, -22.06875977445315))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 279, "", 2.837738513946533, 3.889135059011095, $true, 19.12085205680453, 15, -9.663037743195376, 31, $true, 27, -1.470996706697457, 15.23171699779343, -5.360131765708553, -20.59184876350198))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 361, "", 2.80666732788086, 1.511431856701504, $true, 11.32887560548674, 35, -11.69452457551446, 14, $true, 41, 3.022770271464107, 9.817443748785239, 1.511338414762603, -8.306105334022636))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 362, "", 2.766476392745972, -0.8470847833699174, $false, 12.94624404927052, 34, -10.4116364785305, 13, $true, 40, 4.519468919699076, 13.79332883264043, 5.366553703068993, -8.42677512957144))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 387, "", 2.831804990768433, -4.555893347823535, $false, 10.34061979205806, 9, -6.865022740457345, 52, $true, 15, 2.108246963089976, 14.8965131398816, 6.664140310913512, -8.232372828968089))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 427, "", 2.813739538192749, 0.1926420808429038, $true, 8.564172678250491, 48, -6.267054985453157, 12, $true, 11, -4.201334403350574, 8.371530597407588, -4.393976484193479, -12.76550708160106))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 433, "", 2.986148357391357, -61.0584396702418, $false, 2.296091331119374, 42, -65.1846968836839, 52, $true, 5, -9.73238404977511, 63.35453100136118, 51.32605562046669, -12.02847538089448))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 455, "", 3.021062612533569, -45.4491208645079, $false, 1.113861006607678, 20, -65.58705543488158, 30, $true, 14, -4.455446599409451, 46.56298187111557, 40.99367426509845, -5.569307606017129))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 517, "", 1.848242521286011, 35.24957330867248, $true, 38.76255143003223, 51, -12.08334119783477, 0, $true, 23, 23.08822226158078, 3.512978121359751, -12.1613510470917, -15.67432916845145))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 518, "", 2.033066511154175, 22.5549614244792, $true, 27.14578293557387, 52, -9.090909090909069, 0, $true, 22, 11.8983983087728, 4.590821511094674, -10.65656311570639, -15.24738462680107))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 567, "", 2.524087905883789, -11.35171069763154, $false, 6.751329032658401, 18, -20.2659865852501, 34, $true, 27, -1.485671958039791, 18.10303973028994, 9.866038739591747, -8.237000990698192))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 621, "", 2.451068639755249, -4.696208667123777, $false, 11.44774054863556, 16, -12.45349564669323, 52, $true, 7, 0.696861162856729, 16.14394921575933, 5.393069829980506, -10.75087938577883))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 632, "", 2.596253633499145, -0.2610226368079765, $false, 5.667508501527054, 50, -17.3491801525649, 41, $true, 19, -7.978739474815917, 5.928531138335031, -7.717716838007941, -13.64624797634297))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 686, "", 2.806775331497192, 25.49573925401715, $true, 31.95758318551618, 33, -6.451606948963505, 0, $true, 17, 4.444437271399908, 6.461843931499036, -21.05130198261724, -27.51314591411627))
$null = $newSignalRows.Add(@("twelve_bar_breakout", 715, "", 3.226465225219727, 23.04487268539729, $true, 24.28462571167956, 52, -2.182565318491165, 33, $true, 9, 4.437862564192298, 1.239753026282266, -18.60701012120499, -19.84676314748726))

$startRow = 70
for ($i = 0; $i -lt $newSignalRows.Count; $i++) {
    $row = $newSignalRows[$i]
    $targetRow = $startRow + $i
    for ($col = 0; $col -lt $row.Length; $col++) {
        $wsSignals.Cells.Item($targetRow, $col + 1).Value2 = $row[$col]
    }
}

# ---------------------------------------------------------------------------
# 2) "Summary" sheet: insert the new twelve_bar_breakout row at row 4,
#    pushing random_baseline down to row 5
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Rows.Item(4).Insert()

$summaryRow = @("twelve_bar_breakout", 23, 69.56521739130434, 7.60606683795938, 8.963319693608634, 23.73852399060612, 19.99072970938202, -14.00752976438482, 12.38466287142263, 95.65217391304348, 18.68181818181818, 31.81818181818182, -4.989723397220886, "", "", "", "", "", "", "")
for ($col = 0; $col -lt $summaryRow.Length; $col++) {
    $wsSummary.Cells.Item(4, $col + 1).Value2 = $summaryRow[$col]
}
